$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data, pushing the two
# existing URL rows down to rows 2 and 3.
$ws.Rows("1:1").Insert()

# New header row labelling the exported columns.
$ws.Range("A1").Value = "sites"
$ws.Range("B1").Value = "send"

# New "send" column values next to the (now shifted) site rows.
# "true" is a reserved boolean literal for a plain assignment, so write
# it as a formula first and convert it to a literal value, which keeps
# it stored as text rather than a boolean.
$trueCell = $ws.Range("B2")
$trueCell.Formula = '="true"'
$trueCell.Copy()
$trueCell.PasteSpecial(-4163) # xlPasteValues

$ws.Range("B3").Value = "falsy"

# Match the active selection from the edited workbook.
$ws.Range("B1").Select()
